$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 58
$ws.Cells.Item(5, 9).Value = 58
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 58
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 57
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(19, 8).Value = 2828.8572
$ws.Cells.Item(19, 9).Value = 3000
$ws.Cells.Item(19, 10).Value = 2800.3333
$ws.Cells.Item(19, 11).Value = 3000
$ws.Cells.Item(19, 12).Value = 2800.3333
$ws.Cells.Item(19, 13).Value = -2825
$ws.Cells.Item(19, 14).Value = -3150.3333
$ws.Cells.Item(28, 8).Value = 51488.4
$ws.Cells.Item(28, 9).Value = 91963.63
$ws.Cells.Item(28, 10).Value = 2018.6666
$ws.Cells.Item(28, 11).Value = 91963.63
$ws.Cells.Item(28, 12).Value = 2018.6666
$ws.Cells.Item(28, 13).Value = -91478.63
$ws.Cells.Item(28, 14).Value = -2988.6666
$ws.Cells.Item(34, 8).Value = 20000
$ws.Cells.Item(34, 9).Value = 20000
$ws.Cells.Item(34, 11).Value = 20000
$ws.Cells.Item(34, 13).Value = -19797
$ws.Cells.Item(36, 8).Value = 20000
$ws.Cells.Item(36, 9).Value = 20000
$ws.Cells.Item(36, 11).Value = 20000
$ws.Cells.Item(36, 13).Value = -19285
$ws.Cells.Item(38, 8).Value = 3156.6667
$ws.Cells.Item(38, 9).Value = 2216.2307
$ws.Cells.Item(38, 10).Value = 5601.8
$ws.Cells.Item(38, 11).Value = 6648.6921
$ws.Cells.Item(38, 12).Value = 16805.4
$ws.Cells.Item(38, 13).Value = -6276.6921
$ws.Cells.Item(38, 14).Value = -17549.4
$ws.Cells.Item(47, 8).Value = 15333.333
$ws.Cells.Item(47, 9).Value = 15333.333
$ws.Cells.Item(47, 11).Value = 15333.333
$ws.Cells.Item(47, 13).Value = -14361.333
$ws.Cells.Item(53, 8).Value = 12346386
$ws.Cells.Item(53, 9).Value = 41667360
$ws.Cells.Item(53, 10).Value = 712
$ws.Cells.Item(53, 11).Value = 41667360
$ws.Cells.Item(53, 12).Value = 712
$ws.Cells.Item(53, 13).Value = -41666723
$ws.Cells.Item(53, 14).Value = -1986
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 13).ClearContents()
$ws.Cells.Item(76, 8).Value = 100015370
$ws.Cells.Item(76, 10).Value = 333359680
$ws.Cells.Item(76, 12).Value = 333359680
$ws.Cells.Item(76, 14).Value = -333360310
$ws.Cells.Item(79, 8).Value = 100015370
$ws.Cells.Item(79, 10).Value = 333359680
$ws.Cells.Item(79, 12).Value = 333359680
$ws.Cells.Item(79, 14).Value = -333361864
$ws.Cells.Item(86, 8).Value = 13161394
$ws.Cells.Item(86, 9).Value = 4000
$ws.Cells.Item(86, 10).Value = 17547192
$ws.Cells.Item(86, 11).Value = 4000
$ws.Cells.Item(86, 12).Value = 17547192
$ws.Cells.Item(86, 13).Value = -2877
$ws.Cells.Item(86, 14).Value = -17549438
$ws.Cells.Item(89, 8).Value = 13161394
$ws.Cells.Item(89, 9).Value = 4000
$ws.Cells.Item(89, 10).Value = 17547192
$ws.Cells.Item(89, 11).Value = 20000
$ws.Cells.Item(89, 12).Value = 87735960
$ws.Cells.Item(89, 13).Value = -14384
$ws.Cells.Item(89, 14).Value = -87747192
$ws.Cells.Item(113, 8).Value = 4914.5
$ws.Cells.Item(113, 9).Value = 4899.2
$ws.Cells.Item(113, 10).Value = 4991
$ws.Cells.Item(113, 11).Value = 4899.2
$ws.Cells.Item(113, 12).Value = 4991
$ws.Cells.Item(113, 13).Value = -1645.2
$ws.Cells.Item(113, 14).Value = -11499
$ws.Cells.Item(116, 8).Value = 4058.5715
$ws.Cells.Item(116, 10).Value = 4350.25
$ws.Cells.Item(116, 12).Value = 4350.25
$ws.Cells.Item(116, 14).Value = -11234.25
$ws.Cells.Item(121, 8).Value = 4799.778
$ws.Cells.Item(121, 10).Value = 4799.778
$ws.Cells.Item(121, 12).Value = 14399.334
$ws.Cells.Item(121, 14).Value = -17893.334
$ws.Cells.Item(132, 8).Value = 3092.9048
$ws.Cells.Item(132, 9).Value = 1731.7333
$ws.Cells.Item(132, 11).Value = 5195.199900000001
$ws.Cells.Item(132, 13).Value = -2665.199900000001
$ws.Cells.Item(133, 8).Value = 63891.445
$ws.Cells.Item(133, 10).Value = 63891.445
$ws.Cells.Item(133, 12).Value = 63891.445
$ws.Cells.Item(133, 14).Value = -74011.44500000001
$ws.Cells.Item(134, 8).Value = 74999.8
$ws.Cells.Item(134, 10).Value = 74999.8
$ws.Cells.Item(134, 12).Value = 74999.8
$ws.Cells.Item(134, 14).Value = -85139.8
$ws.Cells.Item(135, 8).Value = 1424.1177
$ws.Cells.Item(135, 9).Value = 654.85187
$ws.Cells.Item(135, 10).Value = 4391.2856
$ws.Cells.Item(135, 11).Value = 5893.66683
$ws.Cells.Item(135, 12).Value = 39521.5704
$ws.Cells.Item(135, 13).Value = -3358.66683
$ws.Cells.Item(135, 14).Value = -44591.5704
$ws.Cells.Item(137, 8).Value = 4257
$ws.Cells.Item(137, 9).Value = 2971.8518
$ws.Cells.Item(137, 10).Value = 5765.6523
$ws.Cells.Item(137, 11).Value = 8915.555399999999
$ws.Cells.Item(137, 12).Value = 17296.9569
$ws.Cells.Item(137, 13).Value = -6365.555399999999
$ws.Cells.Item(137, 14).Value = -22396.9569
$ws.Cells.Item(138, 8).Value = 5146.3335
$ws.Cells.Item(138, 9).Value = 2687.5
$ws.Cells.Item(138, 10).Value = 6215.391
$ws.Cells.Item(138, 11).Value = 8062.5
$ws.Cells.Item(138, 12).Value = 18646.173
$ws.Cells.Item(138, 13).Value = -2922.5
$ws.Cells.Item(138, 14).Value = -28926.173
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 2902.8333
$ws.Cells.Item(141, 9).Value = 1166.75
$ws.Cells.Item(141, 10).Value = 6375
$ws.Cells.Item(141, 11).Value = 3500.25
$ws.Cells.Item(141, 12).Value = 19125
$ws.Cells.Item(141, 13).Value = 1679.75
$ws.Cells.Item(141, 14).Value = -29485

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 4796.5
$ws.Cells.Item(32, 9).Value = 4045.8157
$ws.Cells.Item(32, 11).Value = 4045.8157
$ws.Cells.Item(32, 13).Value = -3758.8157
$ws.Cells.Item(45, 8).Value = 7872.2856
$ws.Cells.Item(45, 9).Value = 4469
$ws.Cells.Item(45, 11).Value = 4469
$ws.Cells.Item(45, 13).Value = -4092
$ws.Cells.Item(61, 8).Value = 3448.7036
$ws.Cells.Item(61, 9).Value = 2339.0476
$ws.Cells.Item(61, 11).Value = 2339.0476
$ws.Cells.Item(61, 13).Value = -2127.0476
$ws.Cells.Item(63, 8).Value = 2680.7778
$ws.Cells.Item(63, 9).Value = 2715.125
$ws.Cells.Item(63, 10).Value = 2406
$ws.Cells.Item(63, 11).Value = 2715.125
$ws.Cells.Item(63, 12).Value = 2406
$ws.Cells.Item(63, 13).Value = -2029.125
$ws.Cells.Item(63, 14).Value = -3778
$ws.Cells.Item(66, 8).Value = 2680.7778
$ws.Cells.Item(66, 9).Value = 2715.125
$ws.Cells.Item(66, 10).Value = 2406
$ws.Cells.Item(66, 11).Value = 13575.625
$ws.Cells.Item(66, 12).Value = 12030
$ws.Cells.Item(66, 13).Value = -10143.625
$ws.Cells.Item(66, 14).Value = -18894
$ws.Cells.Item(74, 8).Value = 1492.4783
$ws.Cells.Item(74, 9).Value = 1348.8334
$ws.Cells.Item(74, 11).Value = 1348.8334
$ws.Cells.Item(74, 13).Value = -474.8334
$ws.Cells.Item(77, 8).Value = 1492.4783
$ws.Cells.Item(77, 9).Value = 1348.8334
$ws.Cells.Item(77, 11).Value = 6744.166999999999
$ws.Cells.Item(77, 13).Value = -2376.166999999999
$ws.Cells.Item(97, 8).Value = 1734.9565
$ws.Cells.Item(97, 9).Value = 1937.0526
$ws.Cells.Item(97, 10).Value = 775
$ws.Cells.Item(97, 11).Value = 1937.0526
$ws.Cells.Item(97, 12).Value = 775
$ws.Cells.Item(97, 13).Value = -1441.0526
$ws.Cells.Item(97, 14).Value = -1767
$ws.Cells.Item(132, 8).Value = 5918.6
$ws.Cells.Item(132, 9).Value = 1779.8182
$ws.Cells.Item(132, 10).Value = 12922.692
$ws.Cells.Item(132, 11).Value = 5339.4546
$ws.Cells.Item(132, 12).Value = 38768.076
$ws.Cells.Item(132, 13).Value = -2809.4546
$ws.Cells.Item(132, 14).Value = -43828.076
$ws.Cells.Item(136, 8).Value = 3448.7036
$ws.Cells.Item(136, 9).Value = 2339.0476
$ws.Cells.Item(136, 11).Value = 7017.1428
$ws.Cells.Item(136, 13).Value = -4467.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 1199
$ws.Cells.Item(86, 9).Value = 1382.7858
$ws.Cells.Item(86, 10).Value = 555.75
$ws.Cells.Item(86, 11).Value = 1382.7858
$ws.Cells.Item(86, 12).Value = 555.75
$ws.Cells.Item(86, 13).Value = -259.7858000000001
$ws.Cells.Item(86, 14).Value = -2801.75
$ws.Cells.Item(89, 8).Value = 1199
$ws.Cells.Item(89, 9).Value = 1382.7858
$ws.Cells.Item(89, 10).Value = 555.75
$ws.Cells.Item(89, 11).Value = 6913.929
$ws.Cells.Item(89, 12).Value = 2778.75
$ws.Cells.Item(89, 13).Value = -1297.929
$ws.Cells.Item(89, 14).Value = -14010.75
$ws.Cells.Item(94, 8).Value = 916.2
$ws.Cells.Item(94, 10).Value = 496
$ws.Cells.Item(94, 12).Value = 496
$ws.Cells.Item(94, 14).Value = -1398
$ws.Cells.Item(107, 8).Value = 1446.3572
$ws.Cells.Item(107, 9).Value = 1359.92
$ws.Cells.Item(107, 11).Value = 1359.92
$ws.Cells.Item(107, 13).Value = 560.0799999999999
$ws.Cells.Item(134, 8).Value = 4503.2705
$ws.Cells.Item(134, 9).Value = 3186.32
$ws.Cells.Item(134, 10).Value = 7246.9165
$ws.Cells.Item(134, 11).Value = 9558.960000000001
$ws.Cells.Item(134, 12).Value = 21740.7495
$ws.Cells.Item(134, 13).Value = -7023.960000000001
$ws.Cells.Item(134, 14).Value = -26810.7495
$ws.Cells.Item(135, 8).Value = 50000
$ws.Cells.Item(135, 10).Value = 50000
$ws.Cells.Item(135, 12).Value = 50000
$ws.Cells.Item(135, 14).Value = -60140
$ws.Cells.Item(137, 8).Value = 75000
$ws.Cells.Item(137, 10).Value = 75000
$ws.Cells.Item(137, 12).Value = 75000
$ws.Cells.Item(137, 14).Value = -85200
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 872.2222
$ws.Cells.Item(7, 9).Value = 831.25
$ws.Cells.Item(7, 10).Value = 1200
$ws.Cells.Item(7, 11).Value = 831.25
$ws.Cells.Item(7, 12).Value = 1200
$ws.Cells.Item(7, 13).Value = -718.25
$ws.Cells.Item(7, 14).Value = -1426
$ws.Cells.Item(31, 8).Value = 2609.878
$ws.Cells.Item(31, 9).Value = 1800.0869
$ws.Cells.Item(31, 10).Value = 3644.611
$ws.Cells.Item(31, 11).Value = 1800.0869
$ws.Cells.Item(31, 12).Value = 3644.611
$ws.Cells.Item(31, 13).Value = -1505.0869
$ws.Cells.Item(31, 14).Value = -4234.611
$ws.Cells.Item(34, 8).Value = 2609.878
$ws.Cells.Item(34, 9).Value = 1800.0869
$ws.Cells.Item(34, 10).Value = 3644.611
$ws.Cells.Item(34, 11).Value = 1800.0869
$ws.Cells.Item(34, 12).Value = 3644.611
$ws.Cells.Item(34, 13).Value = -1598.0869
$ws.Cells.Item(34, 14).Value = -4048.611
$ws.Cells.Item(58, 8).Value = 247282.64
$ws.Cells.Item(58, 9).Value = 501166.56
$ws.Cells.Item(58, 11).Value = 501166.56
$ws.Cells.Item(58, 13).Value = -500963.56
$ws.Cells.Item(62, 8).Value = 5082.1665
$ws.Cells.Item(62, 9).Value = 4497.6665
$ws.Cells.Item(62, 11).Value = 4497.6665
$ws.Cells.Item(62, 13).Value = -3873.6665
$ws.Cells.Item(65, 8).Value = 5082.1665
$ws.Cells.Item(65, 9).Value = 4497.6665
$ws.Cells.Item(65, 11).Value = 22488.3325
$ws.Cells.Item(65, 13).Value = -19368.3325
$ws.Cells.Item(74, 8).Value = 77720
$ws.Cells.Item(74, 10).Value = 77720
$ws.Cells.Item(74, 12).Value = 77720
$ws.Cells.Item(74, 14).Value = -79468
$ws.Cells.Item(77, 8).Value = 77720
$ws.Cells.Item(77, 10).Value = 77720
$ws.Cells.Item(77, 12).Value = 233160
$ws.Cells.Item(77, 14).Value = -241896
$ws.Cells.Item(86, 8).Value = 17925.555
$ws.Cells.Item(86, 9).Value = 7133.25
$ws.Cells.Item(86, 11).Value = 7133.25
$ws.Cells.Item(86, 13).Value = -6010.25
$ws.Cells.Item(88, 8).Value = 21825
$ws.Cells.Item(88, 10).Value = 21825
$ws.Cells.Item(88, 12).Value = 21825
$ws.Cells.Item(88, 14).Value = -22637
$ws.Cells.Item(89, 8).Value = 17925.555
$ws.Cells.Item(89, 9).Value = 7133.25
$ws.Cells.Item(89, 11).Value = 35666.25
$ws.Cells.Item(89, 13).Value = -30050.25
$ws.Cells.Item(91, 8).Value = 21825
$ws.Cells.Item(91, 10).Value = 21825
$ws.Cells.Item(91, 12).Value = 21825
$ws.Cells.Item(91, 14).Value = -24633
$ws.Cells.Item(99, 8).Value = 5348.5293
$ws.Cells.Item(99, 9).Value = 3420.8333
$ws.Cells.Item(99, 10).Value = 6400
$ws.Cells.Item(99, 11).Value = 3420.8333
$ws.Cells.Item(99, 12).Value = 6400
$ws.Cells.Item(99, 13).Value = -1922.8333
$ws.Cells.Item(99, 14).Value = -9396
$ws.Cells.Item(105, 8).Value = 1180.3889
$ws.Cells.Item(105, 9).Value = 1016.4667
$ws.Cells.Item(105, 11).Value = 1016.4667
$ws.Cells.Item(105, 13).Value = 730.5333000000001
$ws.Cells.Item(107, 8).Value = 448.4091
$ws.Cells.Item(107, 9).Value = 302.64285
$ws.Cells.Item(107, 11).Value = 302.64285
$ws.Cells.Item(107, 13).Value = 1617.35715
$ws.Cells.Item(126, 8).Value = 5348.5293
$ws.Cells.Item(126, 9).Value = 3420.8333
$ws.Cells.Item(126, 10).Value = 6400
$ws.Cells.Item(126, 11).Value = 10262.4999
$ws.Cells.Item(126, 12).Value = 19200
$ws.Cells.Item(126, 13).Value = -7792.499899999999
$ws.Cells.Item(126, 14).Value = -24140
$ws.Cells.Item(132, 8).Value = 4512.769
$ws.Cells.Item(132, 9).Value = 2955.3333
$ws.Cells.Item(132, 11).Value = 8865.999899999999
$ws.Cells.Item(132, 13).Value = -6335.999899999999
$ws.Cells.Item(134, 8).Value = 4372.6895
$ws.Cells.Item(134, 9).Value = 3130
$ws.Cells.Item(134, 11).Value = 9390
$ws.Cells.Item(134, 13).Value = -6855
$ws.Cells.Item(136, 8).Value = 247282.64
$ws.Cells.Item(136, 9).Value = 501166.56
$ws.Cells.Item(136, 11).Value = 1503499.68
$ws.Cells.Item(136, 13).Value = -1500949.68
$ws.Cells.Item(139, 8).Value = 95000
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 87.05
$ws.Cells.Item(2, 9).Value = 40.64706
$ws.Cells.Item(2, 11).Value = 243.88236
$ws.Cells.Item(2, 13).Value = -130.88236
$ws.Cells.Item(5, 8).Value = 101511.75
$ws.Cells.Item(5, 9).Value = 101511.75
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 304535.25
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -304423.25
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(34, 8).Value = 55261.523
$ws.Cells.Item(34, 10).Value = 96658.25
$ws.Cells.Item(34, 12).Value = 289974.75
$ws.Cells.Item(34, 14).Value = -290142.75
$ws.Cells.Item(38, 8).Value = 28.1
$ws.Cells.Item(38, 9).Value = 24
$ws.Cells.Item(38, 10).Value = 37.666668
$ws.Cells.Item(38, 11).Value = 72
$ws.Cells.Item(38, 12).Value = 113.000004
$ws.Cells.Item(38, 13).Value = 275
$ws.Cells.Item(38, 14).Value = -807.000004
$ws.Cells.Item(39, 8).Value = 12469
$ws.Cells.Item(39, 9).Value = 1049.5
$ws.Cells.Item(39, 10).Value = 18178.75
$ws.Cells.Item(39, 11).Value = 3148.5
$ws.Cells.Item(39, 12).Value = 54536.25
$ws.Cells.Item(39, 13).Value = -2854.5
$ws.Cells.Item(39, 14).Value = -55124.25
$ws.Cells.Item(122, 8).Value = 126313.625
$ws.Cells.Item(122, 9).Value = 1249.5
$ws.Cells.Item(122, 10).Value = 168001.67
$ws.Cells.Item(122, 11).Value = 11245.5
$ws.Cells.Item(122, 12).Value = 1512015.03
$ws.Cells.Item(122, 13).Value = -8795.5
$ws.Cells.Item(122, 14).Value = -1516915.03
$ws.Cells.Item(132, 8).Value = 6212.5
$ws.Cells.Item(132, 10).Value = 6212.5
$ws.Cells.Item(132, 12).Value = 55912.5
$ws.Cells.Item(132, 14).Value = -60972.5
$ws.Cells.Item(134, 8).Value = 1715.5
$ws.Cells.Item(134, 9).Value = 1715.5
$ws.Cells.Item(134, 11).Value = 5146.5
$ws.Cells.Item(134, 13).Value = -76.5
$ws.Cells.Item(135, 8).Value = 101511.75
$ws.Cells.Item(135, 9).Value = 101511.75
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 913605.75
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 13).Value = -911070.75
$ws.Cells.Item(135, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 3066.5
$ws.Cells.Item(137, 9).Value = 3066.5
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 9199.5
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -4099.5
$ws.Cells.Item(137, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 5951.3335
$ws.Cells.Item(140, 9).Value = 5641.8
$ws.Cells.Item(140, 11).Value = 16925.4
$ws.Cells.Item(140, 13).Value = -11745.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 14950
$ws.Cells.Item(33, 10).Value = 14950
$ws.Cells.Item(33, 12).Value = 14950
$ws.Cells.Item(33, 14).Value = -15454
$ws.Cells.Item(34, 8).Value = 47025.5
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(43, 8).Value = 13212.5
$ws.Cells.Item(43, 9).Value = 6227
$ws.Cells.Item(43, 10).Value = 22526.5
$ws.Cells.Item(43, 11).Value = 6227
$ws.Cells.Item(43, 12).Value = 22526.5
$ws.Cells.Item(43, 13).Value = -6076
$ws.Cells.Item(43, 14).Value = -22828.5
$ws.Cells.Item(48, 8).Value = 18000
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 13).ClearContents()
$ws.Cells.Item(59, 8).Value = 15999
$ws.Cells.Item(59, 10).Value = 15999
$ws.Cells.Item(59, 12).Value = 15999
$ws.Cells.Item(59, 14).Value = -17165
$ws.Cells.Item(70, 8).Value = 52638720
$ws.Cells.Item(70, 10).Value = 66674148
$ws.Cells.Item(70, 12).Value = 66674148
$ws.Cells.Item(70, 14).Value = -66674688
$ws.Cells.Item(73, 8).Value = 52638720
$ws.Cells.Item(73, 10).Value = 66674148
$ws.Cells.Item(73, 12).Value = 66674148
$ws.Cells.Item(73, 14).Value = -66676020
$ws.Cells.Item(76, 8).Value = 47025.5
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 13).ClearContents()
$ws.Cells.Item(79, 8).Value = 47025.5
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 13).ClearContents()
$ws.Cells.Item(80, 8).Value = 4007002
$ws.Cells.Item(80, 9).Value = 2508751
$ws.Cells.Item(80, 11).Value = 2508751
$ws.Cells.Item(80, 13).Value = -2507753
$ws.Cells.Item(83, 8).Value = 4007002
$ws.Cells.Item(83, 9).Value = 2508751
$ws.Cells.Item(83, 11).Value = 12543755
$ws.Cells.Item(83, 13).Value = -12538763

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 78126.07000000001
$ws.Cells.Item(68, 9).Value = 11589.6
$ws.Cells.Item(68, 11).Value = 11589.6
$ws.Cells.Item(68, 13).Value = -10840.6
$ws.Cells.Item(71, 8).Value = 78126.07000000001
$ws.Cells.Item(71, 9).Value = 11589.6
$ws.Cells.Item(71, 11).Value = 57948
$ws.Cells.Item(71, 13).Value = -54204
$ws.Cells.Item(82, 8).Value = 3363.9092
$ws.Cells.Item(82, 9).Value = 2333.1667
$ws.Cells.Item(82, 11).Value = 2333.1667
$ws.Cells.Item(82, 13).Value = -1972.1667
$ws.Cells.Item(85, 8).Value = 3363.9092
$ws.Cells.Item(85, 9).Value = 2333.1667
$ws.Cells.Item(85, 11).Value = 2333.1667
$ws.Cells.Item(85, 13).Value = -1085.1667
$ws.Cells.Item(92, 8).Value = 45000
$ws.Cells.Item(92, 10).Value = 45000
$ws.Cells.Item(92, 12).Value = 45000
$ws.Cells.Item(92, 14).Value = -49992
$ws.Cells.Item(93, 8).Value = 2328.3572
$ws.Cells.Item(93, 9).Value = 2280.3
$ws.Cells.Item(93, 10).Value = 2448.5
$ws.Cells.Item(93, 11).Value = 2280.3
$ws.Cells.Item(93, 12).Value = 2448.5
$ws.Cells.Item(93, 13).Value = -1032.3
$ws.Cells.Item(93, 14).Value = -4944.5
$ws.Cells.Item(132, 8).Value = 3874.3044
$ws.Cells.Item(132, 9).Value = 2403.0908
$ws.Cells.Item(132, 10).Value = 5222.9165
$ws.Cells.Item(132, 11).Value = 7209.2724
$ws.Cells.Item(132, 12).Value = 15668.7495
$ws.Cells.Item(132, 13).Value = -4679.2724
$ws.Cells.Item(132, 14).Value = -20728.7495
$ws.Cells.Item(136, 8).Value = 3529.8235
$ws.Cells.Item(136, 9).Value = 1979.4231
$ws.Cells.Item(136, 11).Value = 5938.2693
$ws.Cells.Item(136, 13).Value = -3388.2693

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 10000000
$ws.Cells.Item(5, 10).Value = 10000000
$ws.Cells.Item(5, 12).Value = 10000000
$ws.Cells.Item(5, 14).Value = -10000224
$ws.Cells.Item(45, 8).Value = 12696.167
$ws.Cells.Item(45, 9).Value = 11999.5
$ws.Cells.Item(45, 10).Value = 12835.5
$ws.Cells.Item(45, 11).Value = 11999.5
$ws.Cells.Item(45, 12).Value = 12835.5
$ws.Cells.Item(45, 13).Value = -11508.5
$ws.Cells.Item(45, 14).Value = -13817.5
$ws.Cells.Item(74, 8).Value = 25000
$ws.Cells.Item(74, 10).Value = 25000
$ws.Cells.Item(74, 12).Value = 25000
$ws.Cells.Item(74, 14).Value = -26872
$ws.Cells.Item(77, 8).Value = 25000
$ws.Cells.Item(77, 10).Value = 25000
$ws.Cells.Item(77, 12).Value = 75000
$ws.Cells.Item(77, 14).Value = -84360
$ws.Cells.Item(81, 8).Value = 8370.941000000001
$ws.Cells.Item(81, 9).Value = 952.0769
$ws.Cells.Item(81, 10).Value = 32482.25
$ws.Cells.Item(81, 11).Value = 1904.1538
$ws.Cells.Item(81, 12).Value = 64964.5
$ws.Cells.Item(81, 13).Value = -843.1538
$ws.Cells.Item(81, 14).Value = -67086.5
$ws.Cells.Item(84, 8).Value = 8370.941000000001
$ws.Cells.Item(84, 9).Value = 952.0769
$ws.Cells.Item(84, 10).Value = 32482.25
$ws.Cells.Item(84, 11).Value = 9520.769
$ws.Cells.Item(84, 12).Value = 324822.5
$ws.Cells.Item(84, 13).Value = -4216.769
$ws.Cells.Item(84, 14).Value = -335430.5
$ws.Cells.Item(86, 8).Value = 99000
$ws.Cells.Item(86, 10).Value = 99000
$ws.Cells.Item(86, 12).Value = 99000
$ws.Cells.Item(86, 14).Value = -101246
$ws.Cells.Item(89, 8).Value = 99000
$ws.Cells.Item(89, 10).Value = 99000
$ws.Cells.Item(89, 12).Value = 495000
$ws.Cells.Item(89, 14).Value = -506232
$ws.Cells.Item(92, 8).Value = 111249.75
$ws.Cells.Item(92, 10).Value = 111249.75
$ws.Cells.Item(92, 12).Value = 111249.75
$ws.Cells.Item(92, 14).Value = -116241.75
$ws.Cells.Item(135, 8).Value = 48999.5
$ws.Cells.Item(135, 10).Value = 48999.5
$ws.Cells.Item(135, 12).Value = 48999.5
$ws.Cells.Item(135, 14).Value = -59139.5
